$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3955935532374564
$ws.Range("D2").Value = 0.4321333824756292
$ws.Range("G2").Value = 0.1242467469831657
$ws.Range("H2").Value = 0.991
